$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    # Force the cell to be stored as text so values that look numeric
    # (e.g. "1.85", "0.0315", "71.77") keep their exact original formatting
    # (trailing zeros, multi-dot thousands separators, etc.) instead of
    # being silently coerced into a true numeric cell by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    # Drop the now-unneeded explicit "Text" number-format style so the
    # cell keeps using the workbook's default style, matching the source.
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.915.96"
$ws.Range("E2").Value = "  +2.09%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.418.72"
$ws.Range("E3").Value = "  +2.31%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.01%  "

# Row 5 - BNB
Set-TextValue "D5" "572.31"
$ws.Range("E5").Value = "  +1.26%  "

# Row 6 - Solana
Set-TextValue "D6" "157.24"
$ws.Range("E6").Value = "  +3.12%  "

# Row 7 - USDC
Set-TextValue "D7" "0.999"
$ws.Range("E7").Value = "  -0.06%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.416.68"
$ws.Range("E8").Value = "  +2.24%  "

# Row 9 - XRP
Set-TextValue "D9" "0.546"
$ws.Range("E9").Value = "  +2.52%  "

# Row 10 - Toncoin
Set-TextValue "D10" "7.36"
$ws.Range("E10").Value = "  -0.90%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +4.14%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -1.27%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.002.49"
$ws.Range("E13").Value = "  +2.20%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -3.39%  "

# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +7.43%  "

# Row 16 - Avalanche
Set-TextValue "D16" "27.23"
$ws.Range("E16").Value = "  +1.44%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "63.894.21"
$ws.Range("E17").Value = "  +2.05%  "

# Row 18 - WrappedEther
Set-TextValue "D18" "3.397.97"
$ws.Range("E18").Value = "  +2.07%  "

# Row 19 - Polkadot
Set-TextValue "D19" "6.28"
$ws.Range("E19").Value = "  -0.93%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +1.77%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "381.74"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22 - Uniswap
Set-TextValue "D22" "8.10"
$ws.Range("E22").Value = "  -3.75%  "

# Row 23 - Dai
Set-TextValue "D23" "1.00"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24 - Litecoin
Set-TextValue "D24" "71.77"
$ws.Range("E24").Value = "  +2.82%  "

# Row 25 - Polygon
Set-TextValue "D25" "0.531"
$ws.Range("E25").Value = "  -0.33%  "

# Row 26 - PEPE
$ws.Range("E26").Value = "  +26.91%  "

# Row 27 - InternetComputer(DFINITY)
Set-TextValue "D27" "9.38"
$ws.Range("E27").Value = "  +1.64%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  -0.22%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.16%  "

# Row 30 - NEARProtocol
$ws.Range("E30").Value = "  +9.38%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "2.02"
$ws.Range("E31").Value = "  +1.42%  "

# Row 32 - Fetch.AI
$ws.Range("E32").Value = "  +3.57%  "

# Row 33 - EthereumClassic
$ws.Range("E33").Value = "  +1.32%  "

# Row 34 - RenderToken
Set-TextValue "D34" "6.39"
$ws.Range("E34").Value = "  -1.01%  "

# Row 35 - USDe
Set-TextValue "D35" "0.998"
$ws.Range("E35").Value = "  -0.02%  "

# Row 36 - Aptos
Set-TextValue "D36" "6.79"
$ws.Range("E36").Value = "  +1.00%  "

# Row 37 - Monero
Set-TextValue "D37" "160.20"
$ws.Range("E37").Value = "  +0.84%  "

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.75%  "

# Row 39 - Maker
Set-TextValue "D39" "2.995.83"
$ws.Range("E39").Value = "  +7.24%  "

# Rows 40/41 - Hedera and Stacks swapped ranking positions
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D40" "1.85"
$ws.Range("E40").Value = "  -1.75%  "

$ws.Range("B41").Value = "Hedera"
$ws.Range("C41").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D41" "0.0757"
$ws.Range("E41").Value = "  +2.85%  "

# Row 42 - EnergySwap
Set-TextValue "D42" "26.99"
$ws.Range("E42").Value = "  +0.11%  "

# Row 43 - VeChain
Set-TextValue "D43" "0.0315"
$ws.Range("E43").Value = "  -1.69%  "

# Row 44 - OKB
Set-TextValue "D44" "42.08"
$ws.Range("E44").Value = "  +3.71%  "

# Row 45 - Mantle
Set-TextValue "D45" "0.760"
$ws.Range("E45").Value = "  +2.35%  "

# Row 46 - Filecoin
Set-TextValue "D46" "4.32"
$ws.Range("E46").Value = "  +1.54%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "23.28"
$ws.Range("E47").Value = "  +5.54%  "

# Row 48 - ONDO
$ws.Range("E48").Value = "  +3.86%  "

# Row 49 - dogwifhat
Set-TextValue "D49" "2.20"
$ws.Range("E49").Value = "  +23.02%  "

# Row 50 - SuiNetwork
Set-TextValue "D50" "0.837"
$ws.Range("E50").Value = "  +4.42%  "

# Row 51 - Cosmos
$ws.Range("E51").Value = "  +0.56%  "
